$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 307 (shifts existing rows 307:394 down to 308:395),
# mirroring a new weekly price observation added above the prior "Camote" entry.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(307, 1).Value = 3
$ws.Cells.Item(307, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(307, 3).Value = "Coquimbo"
$ws.Cells.Item(307, 4).Value = 44551
$ws.Cells.Item(307, 5).Value = 5
$ws.Cells.Item(307, 6).Value = 100112045
$ws.Cells.Item(307, 7).Value = "Zapallo"
$ws.Cells.Item(307, 8).Value = "Camote"
$ws.Cells.Item(307, 9).Value = "1a nueva(o)"
$ws.Cells.Item(307, 10).Value = 250
$ws.Cells.Item(307, 11).Value = 600
$ws.Cells.Item(307, 12).Value = 700
$ws.Cells.Item(307, 13).Value = 648
$ws.Cells.Item(307, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(307, 15).Value = "Provincia de Talca"
$ws.Cells.Item(307, 16).Value = 648
$ws.Cells.Item(307, 17).Value = 1
$ws.Cells.Item(307, 18).Value = "Hortaliza"
